$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.8619958419062073
$ws.Range("C2").Value = 0.105832384099827
$ws.Range("E2").Value = 0.2920739395562038
$ws.Range("F2").Value = 3.258431016905007
$ws.Range("G2").Value = 0.002532622074875974
$ws.Range("J2").Value = 0.1409196760621612
$ws.Range("K2").Value = 0.8902327989225682
$ws.Range("N2").Value = 2.531580154205415
# Row 3
$ws.Range("B3").Value = 0.8165072507206901
$ws.Range("C3").Value = 0.09785692682173419
$ws.Range("E3").Value = 0.2788173909058145
$ws.Range("F3").Value = 3.210343874223824
$ws.Range("G3").Value = 0.002536826843493044
$ws.Range("J3").Value = 0.14025803836698
$ws.Range("K3").Value = 0.8393175211226946
$ws.Range("N3").Value = 2.545131091000101
# Row 4
$ws.Range("B4").Value = 0.7890995301609962
$ws.Range("C4").Value = 0.09302281390478129
$ws.Range("E4").Value = 0.2708622830077374
$ws.Range("F4").Value = 3.182631357288543
$ws.Range("G4").Value = 0.002539542615738815
$ws.Range("J4").Value = 0.1399211948380028
$ws.Range("K4").Value = 0.8085969575795957
$ws.Range("N4").Value = 2.554173099980325
# Row 5
$ws.Range("B5").Value = 0.7780616186934708
$ws.Range("C5").Value = 0.091068531802037
$ws.Range("E5").Value = 0.2676667368863903
$ws.Range("F5").Value = 3.171792867123202
$ws.Range("G5").Value = 0.002540683130291343
$ws.Range("J5").Value = 0.1398013226537103
$ws.Range("K5").Value = 0.7962136996333413
$ws.Range("N5").Value = 2.558038858289763
# Row 6
$ws.Range("B6").Value = 0.7762366793868125
$ws.Range("C6").Value = 0.09074496623696859
$ws.Range("E6").Value = 0.2671389061501372
$ws.Range("F6").Value = 3.17002056077439
$ws.Range("G6").Value = 0.002540874557390055
$ws.Range("J6").Value = 0.1397824668560474
$ws.Range("K6").Value = 0.7941656437500626
$ws.Range("N6").Value = 2.558691688666443
# Row 7
$ws.Range("B7").Value = 0.7889501393365776
$ws.Range("C7").Value = 0.09299639456564535
$ws.Range("E7").Value = 0.2708189998007882
$ws.Range("F7").Value = 3.182483346823773
$ws.Range("G7").Value = 0.002539557860171295
$ws.Range("J7").Value = 0.1399195078465283
$ws.Range("K7").Value = 0.8084294042617444
$ws.Range("N7").Value = 2.554224502375668
# Row 8
$ws.Range("B8").Value = 0.8462026684386217
$ws.Range("C8").Value = 0.1030693109198495
$ws.Range("E8").Value = 0.2874646930246527
$ws.Range("F8").Value = 3.241473400947797
$ws.Range("G8").Value = 0.002534044129823617
$ws.Range("J8").Value = 0.1406771028086027
$ws.Range("K8").Value = 0.8725644310184748
$ws.Range("N8").Value = 2.536102487845156
# Row 9
$ws.Range("B9").Value = 0.9626462717680795
$ws.Range("C9").Value = 0.1233286924992285
$ws.Range("E9").Value = 0.3215799467233325
$ws.Range("F9").Value = 3.37161118914176
$ws.Range("G9").Value = 0.002524290002663905
$ws.Range("J9").Value = 0.1427163068947621
$ws.Range("K9").Value = 1.00266433788974
$ws.Range("N9").Value = 2.506309850059296
# Row 10
$ws.Range("B10").Value = 1.050786989470794
$ws.Range("C10").Value = 0.1385343593713912
$ws.Range("E10").Value = 0.3475587076929116
$ws.Range("F10").Value = 3.476152258189529
$ws.Range("G10").Value = 0.002517761549205622
$ws.Range("J10").Value = 0.1445563052047518
$ws.Range("K10").Value = 1.100947357154951
$ws.Range("N10").Value = 2.487948807996432
# Row 11
$ws.Range("B11").Value = 1.091457593848475
$ws.Range("C11").Value = 0.1455243310665253
$ws.Range("E11").Value = 0.3595794031715656
$ws.Range("F11").Value = 3.525675502918887
$ws.Range("G11").Value = 0.002514928552698257
$ws.Range("J11").Value = 0.1454685529526074
$ws.Range("K11").Value = 1.146258215514365
$ws.Range("N11").Value = 2.480367512488641
# Row 12
$ws.Range("B12").Value = 1.106941846726215
$ws.Range("C12").Value = 0.14818192100293
$ws.Range("E12").Value = 0.3641607409251151
$ws.Range("F12").Value = 3.544713369162224
$ws.Range("G12").Value = 0.002513875328645407
$ws.Range("J12").Value = 0.1458248879365343
$ws.Range("K12").Value = 1.163503609565538
$ws.Range("N12").Value = 2.477608059067734
# Row 13
$ws.Range("B13").Value = 1.10360332837331
$ws.Range("C13").Value = 0.1476090856105259
$ws.Range("E13").Value = 0.3631727588624187
$ws.Range("F13").Value = 3.540600544964519
$ws.Range("G13").Value = 0.002514101290544868
$ws.Range("J13").Value = 0.1457476595224279
$ws.Range("K13").Value = 1.159785618291636
$ws.Range("N13").Value = 2.478197393062644
# Row 14
$ws.Range("B14").Value = 1.092729822737851
$ws.Range("C14").Value = 0.1457427584051914
$ws.Range("E14").Value = 0.3599557230077437
$ws.Range("F14").Value = 3.527236048586133
$ws.Range("G14").Value = 0.002514841511665747
$ws.Range("J14").Value = 0.1454976502615253
$ws.Range("K14").Value = 1.147675253840788
$ws.Range("N14").Value = 2.480138254619661
# Row 15
$ws.Range("B15").Value = 1.086080333815858
$ws.Range("C15").Value = 0.1446009696152544
$ws.Range("E15").Value = 0.3579890252313049
$ws.Range("F15").Value = 3.519087006505117
$ws.Range("G15").Value = 0.002515297463947581
$ws.Range("J15").Value = 0.1453459322020692
$ws.Range("K15").Value = 1.140268676313354
$ws.Range("N15").Value = 2.481341613793603
# Row 16
$ws.Range("B16").Value = 1.04814068181031
$ws.Range("C16").Value = 0.1380790280396411
$ws.Range("E16").Value = 0.3467772249262993
$ws.Range("F16").Value = 3.472955517883719
$ws.Range("G16").Value = 0.002517949436334909
$ws.Range("J16").Value = 0.1444982072779979
$ws.Range("K16").Value = 1.097998345230565
$ws.Range("N16").Value = 2.488459826653411
# Row 17
$ws.Range("B17").Value = 1.025013545088996
$ws.Range("C17").Value = 0.1340968004500382
$ws.Range("E17").Value = 0.3399512346209406
$ws.Range("F17").Value = 3.44516029329327
$ws.Range("G17").Value = 0.00251961130581659
$ws.Range("J17").Value = 0.1439974711731011
$ws.Range("K17").Value = 1.072221378068946
$ws.Range("N17").Value = 2.493024516329797
# Row 18
$ws.Range("B18").Value = 1.011765514092986
$ws.Range("C18").Value = 0.1318131742337982
$ws.Range("E18").Value = 0.336044184603125
$ws.Range("F18").Value = 3.429358257948877
$ws.Range("G18").Value = 0.002520580054358943
$ws.Range("J18").Value = 0.1437165371691975
$ws.Range("K18").Value = 1.057451715817081
$ws.Range("N18").Value = 2.495722566096745
# Row 19
$ws.Range("B19").Value = 1.007289232749429
$ws.Range("C19").Value = 0.1310411484058136
$ws.Range("E19").Value = 0.3347245955847171
$ws.Range("F19").Value = 3.424039692078594
$ws.Range("G19").Value = 0.002520910272309624
$ws.Range("J19").Value = 0.1436226309757203
$ws.Range("K19").Value = 1.052460657224714
$ws.Range("N19").Value = 2.496648526184288
# Row 20
$ws.Range("B20").Value = 1.027469866958711
$ws.Range("C20").Value = 0.1345200056727265
$ws.Range("E20").Value = 0.3406758969312875
$ws.Range("F20").Value = 3.448099979462199
$ws.Range("G20").Value = 0.002519433064038676
$ws.Range("J20").Value = 0.1440500425294502
$ws.Range("K20").Value = 1.074959520888825
$ws.Range("N20").Value = 2.492531083862076
# Row 21
$ws.Range("B21").Value = 1.095921373455269
$ws.Range("C21").Value = 0.1462906538651509
$ws.Range("E21").Value = 0.3608998459858412
$ws.Range("F21").Value = 3.531153791012684
$ws.Range("G21").Value = 0.002514623560392905
$ws.Range("J21").Value = 0.1455707880250472
$ws.Range("K21").Value = 1.151229990644481
$ws.Range("N21").Value = 2.479565148514283
# Row 22
$ws.Range("B22").Value = 1.141143453836605
$ws.Range("C22").Value = 0.1540455458425356
$ws.Range("E22").Value = 0.3742885997148733
$ws.Range("F22").Value = 3.587093253144246
$ws.Range("G22").Value = 0.002511594294185286
$ws.Range("J22").Value = 0.146628171260403
$ws.Range("K22").Value = 1.201585481545862
$ws.Range("N22").Value = 2.471740837877746
# Row 23
$ws.Range("B23").Value = 1.116963026609881
$ws.Range("C23").Value = 0.1499008772069317
$ws.Range("E23").Value = 0.3671270362409444
$ws.Range("F23").Value = 3.557084962903531
$ws.Range("G23").Value = 0.002513200672483053
$ws.Range("J23").Value = 0.1460579935377027
$ws.Range("K23").Value = 1.174663070533569
$ws.Range("N23").Value = 2.475857201825136
# Row 24
$ws.Range("B24").Value = 1.026359214155093
$ws.Range("C24").Value = 0.1343286565094388
$ws.Range("E24").Value = 0.3403482232227191
$ws.Range("F24").Value = 3.446770393605306
$ws.Range("G24").Value = 0.002519513605508228
$ws.Range("J24").Value = 0.1440262533900381
$ws.Range("K24").Value = 1.073721451267176
$ws.Range("N24").Value = 2.492753934985728
# Row 25
$ws.Range("B25").Value = 0.9306938311340502
$ws.Range("C25").Value = 0.1177925071876018
$ws.Range("E25").Value = 0.3121914986108578
$ws.Range("F25").Value = 3.334845842189537
$ws.Range("G25").Value = 0.002526816208347819
$ws.Range("J25").Value = 0.1421049412351039
$ws.Range("K25").Value = 0.9669992080443421
$ws.Range("N25").Value = 2.513751934095879
